$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-6 with corrected financial figures
$ws.Range("D2").Value = 2444
$ws.Range("E2").Value = -84
$ws.Range("F2").Value = -84
$ws.Range("G2").Value = -170
$ws.Range("H2").Value = -115
$ws.Range("I2").Value = -115
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 3644
$ws.Range("L2").Value = 1931
$ws.Range("M2").Value = 1713
$ws.Range("N2").Value = 1713
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 53
$ws.Range("Q2").Value = 187
$ws.Range("R2").Value = -33
$ws.Range("S2").Value = -133
$ws.Range("T2").Value = 48
$ws.Range("U2").Value = 139
$ws.Range("V2").Value = 1394
$ws.Range("W2").Value = -3.44
$ws.Range("X2").Value = -4.72
$ws.Range("Y2").Value = -9.45
$ws.Range("Z2").Value = -3.72
$ws.Range("AA2").Value = 112.72
$ws.Range("AB2").Value = 1030.65
$ws.Range("AC2").Value = -2168
$ws.Range("AD2").Value = -2.86
$ws.Range("AE2").Value = 32326
$ws.Range("AF2").Value = 0.19
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 5300000

$ws.Range("D3").Value = 2444
$ws.Range("E3").Value = 43
$ws.Range("F3").Value = 43
$ws.Range("G3").Value = -76
$ws.Range("H3").Value = -71
$ws.Range("I3").Value = -70
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3705
$ws.Range("L3").Value = 2068
$ws.Range("M3").Value = 1638
$ws.Range("N3").Value = 1638
$ws.Range("O3").Value = -1
$ws.Range("P3").Value = 53
$ws.Range("Q3").Value = -92
$ws.Range("R3").Value = -12
$ws.Range("S3").Value = 101
$ws.Range("T3").Value = 6
$ws.Range("U3").Value = -98
$ws.Range("V3").Value = 1498
$ws.Range("W3").Value = 1.76
$ws.Range("X3").Value = -2.89
$ws.Range("Y3").Value = -4.19
$ws.Range("Z3").Value = -1.92
$ws.Range("AA3").Value = 126.26
$ws.Range("AB3").Value = 894.86
$ws.Range("AC3").Value = -1325
$ws.Range("AD3").Value = -9.67
$ws.Range("AE3").Value = 30907
$ws.Range("AF3").Value = 0.41
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 5300000

$ws.Range("D4").Value = 2456
$ws.Range("E4").Value = 12
$ws.Range("F4").Value = 12
$ws.Range("G4").Value = -40
$ws.Range("H4").Value = -41
$ws.Range("I4").Value = -41
$ws.Range("J4").Value = -1
$ws.Range("K4").Value = 3716
$ws.Range("L4").Value = 1926
$ws.Range("M4").Value = 1790
$ws.Range("N4").Value = 1791
$ws.Range("O4").Value = -1
$ws.Range("P4").Value = 53
$ws.Range("Q4").Value = -79
$ws.Range("R4").Value = 162
$ws.Range("S4").Value = -87
$ws.Range("T4").Value = 9
$ws.Range("U4").Value = -88
$ws.Range("V4").Value = 1411
$ws.Range("W4").Value = 0.49
$ws.Range("X4").Value = -1.69
$ws.Range("Y4").Value = -2.38
$ws.Range("Z4").Value = -1.12
$ws.Range("AA4").Value = 107.58
$ws.Range("AB4").Value = 818.79
$ws.Range("AC4").Value = -770
$ws.Range("AD4").Value = -23.32
$ws.Range("AE4").Value = 33801
$ws.Range("AF4").Value = 0.53
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 5300000

$ws.Range("D5").Value = 2086
$ws.Range("E5").Value = -99
$ws.Range("F5").Value = -99
$ws.Range("G5").Value = -136
$ws.Range("H5").Value = -96
$ws.Range("I5").Value = -96
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3560
$ws.Range("L5").Value = 1876
$ws.Range("M5").Value = 1683
$ws.Range("N5").Value = 1685
$ws.Range("O5").Value = -1
$ws.Range("P5").Value = 53
$ws.Range("Q5").Value = 31
$ws.Range("R5").Value = 196
$ws.Range("S5").Value = -246
$ws.Range("T5").Value = 13
$ws.Range("U5").Value = 18
$ws.Range("V5").Value = 1166
$ws.Range("W5").Value = -4.76
$ws.Range("X5").Value = -4.6
$ws.Range("Y5").Value = -5.51
$ws.Range("Z5").Value = -2.64
$ws.Range("AA5").Value = 111.44
$ws.Range("AB5").Value = 644.61
$ws.Range("AC5").Value = -1808
$ws.Range("AD5").Value = -8.23
$ws.Range("AE5").Value = 31785
$ws.Range("AF5").Value = 0.47
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 5300000

$ws.Range("D6").Value = 2091
$ws.Range("E6").Value = -134
$ws.Range("F6").Value = -134
$ws.Range("G6").Value = 103
$ws.Range("H6").Value = 124
$ws.Range("I6").Value = 125
$ws.Range("K6").Value = 2675
$ws.Range("L6").Value = 907
$ws.Range("M6").Value = 1768
$ws.Range("N6").Value = 1770
$ws.Range("P6").Value = 53
$ws.Range("Q6").Value = -187
$ws.Range("R6").Value = 1712
$ws.Range("S6").Value = -698
$ws.Range("T6").Value = 17
$ws.Range("U6").Value = -204
$ws.Range("V6").Value = 492
$ws.Range("W6").Value = -6.42
$ws.Range("X6").Value = 5.95
$ws.Range("Y6").Value = 7.26
$ws.Range("Z6").Value = 3.99
$ws.Range("AA6").Value = 51.28
$ws.Range("AB6").Value = 2762.97
$ws.Range("AC6").Value = 2365
$ws.Range("AD6").Value = 6.17
$ws.Range("AE6").Value = 34528
$ws.Range("AF6").Value = 0.42
$ws.Range("AG6").Value = 300
$ws.Range("AH6").Value = 2.05
$ws.Range("AI6").Value = 12.28
$ws.Range("AJ6").Value = 5300000

# Rows 7-9 (years with no reliable data) - clear all value columns, keep A/B/C labels
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
